$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows 27-30 duplicate existing employee records (row26, row24, row25, row26)
# matching the source edit that appended these records to the bottom of the table.
$newRows = @(
    @(919133, "Mr.", "Julio",   "L", "Rossignol", "M", "julio.rossignol@shaw.ca",    "Clyde Rossignol", "Amiee Rossignol", "Sirois",    27993),
    @(925382, "Mr.", "Mikel",   "P", "Rolfe",     "M", "mikel.rolfe@gmail.com",      "Dwayne Rolfe",    "Robbin Rolfe",    "Whitehill", 23270),
    @(913180, "Ms.", "Paulina", "U", "Mcgeorge",  "F", "paulina.mcgeorge@yahoo.com", "Joe Mcgeorge",    "Barbar Mcgeorge", "Rodger",    24313),
    @(919133, "Mr.", "Julio",   "L", "Rossignol", "M", "julio.rossignol@shaw.ca",    "Clyde Rossignol", "Amiee Rossignol", "Sirois",    27993)
)

$startRow = 27
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value2 = $rowData[$c]
    }
    # Keep the DateofBirth column formatted the same way as the rest of column K.
    $ws.Cells.Item($r, 11).NumberFormat = $ws.Cells.Item(26, 11).NumberFormat
}

# Move the selection to the row just below the newly appended data (mirrors the
# original author selecting the next empty row after finishing the paste).
[void]$ws.Rows.Item(31).Select()
